$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.521.71"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +4.28%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.484.18"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +2.64%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.74"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.66"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +4.26%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.483.43"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +2.57%  "
$ws.Range("E9").Value = "  +7.41%  "
$ws.Range("E10").Value = "  +0.23%  "
$ws.Range("E11").Value = "  +6.40%  "
$ws.Range("E12").Value = "  +3.78%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.088.55"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.64%  "
$ws.Range("E14").Value = "  +0.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.12"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +4.45%  "
$ws.Range("E16").Value = "  +3.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.567.81"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +4.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.464.26"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +3.50%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.31"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +3.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.07"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +3.71%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "391.79"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +4.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.91"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.91%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "73.22"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +4.27%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.998"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.19%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.534"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +4.20%  "
$ws.Range("E26").Value = "  +5.50%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.12"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +7.46%  "
$ws.Range("E28").Value = "  +2.04%  "
$ws.Range("E29").Value = "  +0.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.32"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +4.12%  "
$ws.Range("E31").Value = "  +6.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.05"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.92%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.57"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +3.35%  "
$ws.Range("E34").Value = "  +4.86%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").ClearFormats()
$ws.Range("E36").Value = "  +9.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "161.91"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.40%  "
$ws.Range("E38").Value = "  +3.04%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.91"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +6.19%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.74"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +4.53%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.65"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +6.40%  "
$ws.Range("B42").Value = "Hedera"
$ws.Range("C42").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0743"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +3.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.46"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +3.38%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "26.98"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +3.58%  "
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "43.06"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.69%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.775.11"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0311"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.80%  "
$ws.Range("B48").Value = "dogwifhat"
$ws.Range("C48").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.47"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.15%  "
$ws.Range("B49").Value = "Bittensor"
$ws.Range("C49").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("E49").Value = "  +5.81%  "
$ws.Range("E50").Value = "  +5.35%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.96"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +13.75%  "
